# "Hash Map added for signup page test"
# Rename Sheet2 -> Signup, populate it with signup-form sample data
# (a HashMap of field -> value used by the UI test), and make it the
# active sheet/tab.

$wb = $excel.ActiveWorkbook

$credentials = $wb.Worksheets.Item(1)
$signup      = $wb.Worksheets.Item(2)

$signup.Name = "Signup"

# Header row - keys of the test's hash map, written left to right.
$signup.Range("A1").Value = "Email"
$signup.Range("B1").Value = "Gender"
$signup.Range("C1").Value = "FirstName"
$signup.Range("D1").Value = "LastName"
$signup.Range("E1").Value = "Password"
$signup.Range("F1").Value = "Day"
$signup.Range("G1").Value = "Month"
$signup.Range("H1").Value = "Year"

# Data row - corresponding values (filled in order, email/hyperlink cell last).
$signup.Range("B2").Value = "Mr"
$signup.Range("C2").Value = "Minuga"
$signup.Range("D2").Value = "Lakvindu"
$signup.Range("E2").Value = "Minuga@123"
$signup.Range("F2").Value = 30
$signup.Range("G2").Value = "September"
$signup.Range("H2").Value = 2005
$signup.Range("A2").Value = "minuga@abf.com"

# Excel auto-hyperlinked the two "@"-containing values when they were typed in.
$signup.Hyperlinks.Add($signup.Range("A2"), "mailto:minuga@abf.com")
$signup.Hyperlinks.Add($signup.Range("E2"), "mailto:Minuga@123")
$signup.Range("A2").Style = "Hyperlink"
$signup.Range("E2").Style = "Hyperlink"

# Best-fit the Month column like the rest of the sheet.
$signup.Columns.Item(7).AutoFit() | Out-Null

# Switch focus to the new Signup tab and land the selection like the author left it.
$signup.Activate() | Out-Null
$signup.Range("N18").Select() | Out-Null

Write-Output "done"
